$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.145.28"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "1.867.86"
$ws.Range("E3").Value = "  -1.93%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.74"
$ws.Range("E5").Value = "  -1.99%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("E7").Value = "  -1.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3755"
$ws.Range("E8").Value = "  -0.78%  "

$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8912"
$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.77"
$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.870.04"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07536"
$ws.Range("E13").Value = "  -1.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.301"
$ws.Range("E14").Value = "  -2.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.62"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008487"
$ws.Range("E17").Value = "  -2.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.09"
$ws.Range("E18").Value = "  -3.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "27.182.64"
$ws.Range("E20").Value = "  -2.13%  "

$ws.Range("E21").Value = "  -2.77%  "

$ws.Range("D22").Value = "2.102.08"
$ws.Range("E22").Value = "  -2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.46"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.462"
$ws.Range("E24").Value = "  -2.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.837"
$ws.Range("E25").Value = "  -1.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.24"
$ws.Range("E26").Value = "  -4.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("E27").Value = "  -1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.083"
$ws.Range("E28").Value = "  -3.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.91"
$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.660"
$ws.Range("E30").Value = "  -4.03%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.687"
$ws.Range("E31").Value = "  -3.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09234"
$ws.Range("E32").Value = "  +2.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05125"
$ws.Range("E33").Value = "  -3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.081"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.161"
$ws.Range("E35").Value = "  -5.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7252"
$ws.Range("E36").Value = "  -7.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02031"
$ws.Range("E37").Value = "  -3.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.103"
$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.514"
$ws.Range("E39").Value = "  -3.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.077"
$ws.Range("E40").Value = "  -1.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5285"
$ws.Range("E41").Value = "  -4.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.505"
$ws.Range("E42").Value = "  -3.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.44"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.306"
$ws.Range("E44").Value = "  -2.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1469"
$ws.Range("E45").Value = "  -3.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9990"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("E47").Value = "  -4.14%  "

$ws.Range("E48").Value = "  -5.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.558"
$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.67"
$ws.Range("E50").Value = "  -0.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.63"
$ws.Range("E51").Value = "  -4.72%  "
